$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to remain plain text,
# matching the workbook author convention of storing prices as inline strings
# (e.g. "235.23", "1.863.65") rather than numeric values. (Applied as separate
# contiguous ranges since comma-separated union ranges only format the first area.)
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D8:D17").NumberFormat = "@"
$ws.Range("D19:D23").NumberFormat = "@"
$ws.Range("D25:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.237.83'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.864.69'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '235.23'
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").Value = '0.2835'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.06516'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '21.47'
$ws.Range("E10").Value = '  +3.32%  '
$ws.Range("D11").Value = '0.07855'
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").Value = '97.40'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = '1.866.25'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").Value = '5.093'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '0.6722'
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").Value = '280.43'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").Value = '30.233.64'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '5.482'
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("D20").Value = '12.69'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '2.116.05'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").Value = '0.000007277'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").Value = '9.172'
$ws.Range("E25").Value = '  -2.47%  '
$ws.Range("D26").Value = '164.23'
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("D27").Value = '19.10'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("D28").Value = '1.928'
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("D29").Value = '1.378'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '0.09678'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").Value = '4.403'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '1.477'
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").Value = '4.092'
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("D34").Value = '0.04685'
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").Value = '1.116'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = '0.7053'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").Value = '2.727'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").Value = '0.01849'
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("D39").Value = '2.533'
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").Value = '6.243'
$ws.Range("E40").Value = '  -6.61%  '
$ws.Range("D41").Value = '73.07'
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("D42").Value = '1.940'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("D43").Value = '0.8470'
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4169'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '103.91'
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = '0.9996'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").Value = '7.191'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").Value = '9.163'
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("D49").Value = '935.80'
$ws.Range("E49").Value = '  -5.69%  '
$ws.Range("D50").Value = '34.13'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '0.1124'
$ws.Range("E51").Value = '  -2.28%  '
